$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.255.17"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.037.28"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.397"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0814"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.863"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.335.48"
$ws.Range("E14").Value = "  +3.80%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("D17").Value = "2.036.86"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "37.182.78"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "0.0₃0871"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.90%  "
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("E30").Value = "  +6.17%  "
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0667"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0216"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("E46").Value = "  +4.20%  "
$ws.Range("D47").Value = "1.393.33"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +18.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.21%  "
